$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row 47 by copying the formatting of row 46 (keeps same
# cell styles as the rest of the table: s="3"/"4"/"5").
$ws.Rows(46).Copy() | Out-Null
$ws.Rows(47).Insert(-4121) | Out-Null

# Fill in the values for Post 37. Set them in the same order the
# original author's shared-strings table grew (dev.to link, then
# title, then hashnode link) so new shared string indices line up.
$ws.Range("F47").Value = "https://dev.to/rahulmishra05/critical-section-problem-operating-system-mo03-p04-4fhg"
$ws.Range("C47").Value = "Critical Section Problem | Operating System - M03 P04"
$ws.Range("E47").Value = "https://programmingport.hashnode.dev/critical-section-problem-or-operating-system-m03-p04"
$ws.Range("B47").Value = 37
$ws.Range("D47").Value = Get-Date -Year 2020 -Month 12 -Day 2 -Hour 0 -Minute 0 -Second 0

# Grow the table / autofilter range to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B10:F47"))

# Update the view so the new row is the active selection.
$ws.Activate() | Out-Null
$ws.Range("E47").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 4
